# "corrected data cleaning for pre/post/total fixation data"
#
# 1) Header row (row 1) loses the bold/bordered/centered "style 1" formatting
#    that was applied to A1:AH1 - cells fall back to the default style (0).
# 2) A1's text ("Unnamed: 0") is cleared entirely (empty cell).
# 3) A batch of summary statistics (pre/post/total fixation columns H, I, M,
#    R, AB, plus a few percentage columns in row 6) are corrected to new
#    values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1 & 2: strip header formatting, clear A1's text -----------------------
$ws.Range("A1:AH1").ClearFormats()
$ws.Range("A1").ClearContents()

# --- 3: corrected numeric values --------------------------------------------

# Row 3 - Revisit count
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 33
$ws.Range("M3").Value = 34
$ws.Range("R3").Value = 32
$ws.Range("AB3").Value = 37

# Row 4 - Fixation count
$ws.Range("H4").Value = 110
$ws.Range("I4").Value = 96
$ws.Range("M4").Value = 74
$ws.Range("R4").Value = 91
$ws.Range("AB4").Value = 440

# Row 5 - Dwell time (ms)
$ws.Range("H5").Value = 24111.33
$ws.Range("I5").Value = 21925.71
$ws.Range("M5").Value = 16753.3
$ws.Range("R5").Value = 21424.83
$ws.Range("AB5").Value = 132958.78

# Row 6 - Dwell time (%)
$ws.Range("F6").Value = 0.39
$ws.Range("H6").Value = 9.51
$ws.Range("I6").Value = 8.65
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0.08
$ws.Range("L6").Value = 0.57
$ws.Range("M6").Value = 6.61
$ws.Range("N6").Value = 0.08
$ws.Range("P6").Value = 0.28
$ws.Range("R6").Value = 8.449999999999999
$ws.Range("T6").Value = 0.08
$ws.Range("W6").Value = 0.51
$ws.Range("X6").Value = 0.31
$ws.Range("AB6").Value = 52.42
$ws.Range("AH6").Value = 0.11

# Row 7 - Fixation duration (ms)
$ws.Range("H7").Value = 219.19
$ws.Range("I7").Value = 228.39
$ws.Range("M7").Value = 226.4
$ws.Range("R7").Value = 235.44
$ws.Range("AB7").Value = 302.18
